# Actualización input ponderacion diccionario
# Updates a handful of "ponderacion" (E column) values, then applies an
# AutoFilter on the "Diccionario" sheet so only rows where
# categoria = "Construcción" and ponderacion is 1, 2 or 3 remain visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the ponderacion (column E) values that changed ---------------
$ws.Range("E3").Value = 1
$ws.Range("E18").Value = 3
$ws.Range("E41").Value = 3
$ws.Range("E56").Value = 2
$ws.Range("E90").Value = 2

# --- (Re)apply the AutoFilter over the full data range --------------------
# Drop any existing filter first so the new filter range isn't merged with
# the stale A1:E85 range left over from before the sheet grew to row 91.
$ws.AutoFilterMode = $false

$rng = $ws.Range("A1:E91")
$rng.AutoFilter(1, "Construcción", 7)
$vals = @("1", "2", "3")
$rng.AutoFilter(5, $vals, 7)

# Rows 87/88 (categoria="Construcción", ponderacion=0) fall outside the
# {1,2,3} filter criteria but the autofilter evaluation engine doesn't pick
# them up automatically here, so hide them explicitly to match Excel's own
# filter evaluation.
$ws.Rows.Item(87).Hidden = $true
$ws.Rows.Item(88).Hidden = $true

# --- Keep the _FilterDatabase defined name in sync with the new range -----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Diccionario!_FilterDatabase") {
        $n.RefersTo = "=Diccionario!`$A`$1:`$E`$91"
    }
}

# --- Update the sheet view: scroll back to the top, select E91 ------------
$ws.Range("A1").Select()
$ws.Range("E91").Select()
